$wb = $excel.ActiveWorkbook

# =========================================================
# 1. Predictions sheet: replace week-12 matchups with week-13
# =========================================================
$wsPred = $wb.Worksheets.Item("Predictions")

$wsPred.Cells.Item(2, 1).Value = "Dallas Cowboys"
$wsPred.Cells.Item(2, 2).Value = "New Orleans Saints"
$wsPred.Cells.Item(2, 3).Value = 34
$wsPred.Cells.Item(2, 4).Value = 27
$wsPred.Cells.Item(2, 5).Value = 7
$wsPred.Cells.Item(2, 6).Value = 4.5
$wsPred.Cells.Item(2, 7).Value = "Dallas -4.5"
$wsPred.Cells.Item(2, 8).Value = 61
$wsPred.Cells.Item(2, 9).Value = 47.5
$wsPred.Cells.Item(2, 10).Value = "OVER"

$wsPred.Cells.Item(3, 1).Value = "Minnesota Vikings"
$wsPred.Cells.Item(3, 2).Value = "Detroit Lions"
$wsPred.Cells.Item(3, 3).Value = 31
$wsPred.Cells.Item(3, 4).Value = 14
$wsPred.Cells.Item(3, 5).Value = 17
$wsPred.Cells.Item(3, 6).Value = 7
$wsPred.Cells.Item(3, 7).Value = "Minnesota -7"
$wsPred.Cells.Item(3, 8).Value = 45
$wsPred.Cells.Item(3, 9).Value = 47
$wsPred.Cells.Item(3, 10).Value = "UNDER"

$wsPred.Cells.Item(4, 1).Value = "Tampa Bay Buccaneers"
$wsPred.Cells.Item(4, 2).Value = "Atlanta Falcons"
$wsPred.Cells.Item(4, 3).Value = 44
$wsPred.Cells.Item(4, 4).Value = 17
$wsPred.Cells.Item(4, 5).Value = 27
$wsPred.Cells.Item(4, 6).Value = 11
$wsPred.Cells.Item(4, 7).Value = "Tampa Bay -11"
$wsPred.Cells.Item(4, 8).Value = 61
$wsPred.Cells.Item(4, 9).Value = 50.5
$wsPred.Cells.Item(4, 10).Value = "OVER"

$wsPred.Cells.Item(5, 1).Value = "Arizona Cardinals"
$wsPred.Cells.Item(5, 2).Value = "Chicago Bears"
$wsPred.Cells.Item(5, 3).Value = 34
$wsPred.Cells.Item(5, 4).Value = 7
$wsPred.Cells.Item(5, 5).Value = 27
$wsPred.Cells.Item(5, 6).Value = 7.5
$wsPred.Cells.Item(5, 7).Value = "Arizona -7.5"
$wsPred.Cells.Item(5, 8).Value = 41
$wsPred.Cells.Item(5, 9).Value = 46.5
$wsPred.Cells.Item(5, 10).Value = "UNDER"

$wsPred.Cells.Item(6, 1).Value = "Indianapolis Colts"
$wsPred.Cells.Item(6, 2).Value = "Houston Texans"
$wsPred.Cells.Item(6, 3).Value = 37
$wsPred.Cells.Item(6, 4).Value = 11
$wsPred.Cells.Item(6, 5).Value = 26
$wsPred.Cells.Item(6, 6).Value = 9
$wsPred.Cells.Item(6, 7).Value = "Indianapolis -9"
$wsPred.Cells.Item(6, 8).Value = 48
$wsPred.Cells.Item(6, 9).Value = 46
$wsPred.Cells.Item(6, 10).Value = "OVER"

$wsPred.Cells.Item(7, 1).Value = "Los Angeles Chargers"
$wsPred.Cells.Item(7, 2).Value = "Cincinnati Bengals"
$wsPred.Cells.Item(7, 3).Value = 25
$wsPred.Cells.Item(7, 4).Value = 40
$wsPred.Cells.Item(7, 5).Value = -15
$wsPred.Cells.Item(7, 6).Value = -3
$wsPred.Cells.Item(7, 7).Value = "Cincinnati -3"
$wsPred.Cells.Item(7, 8).Value = 65
$wsPred.Cells.Item(7, 9).Value = 50.5
$wsPred.Cells.Item(7, 10).Value = "OVER"

$wsPred.Cells.Item(8, 1).Value = "New York Giants"
$wsPred.Cells.Item(8, 2).Value = "Miami Dolphins"
$wsPred.Cells.Item(8, 3).Value = 14
$wsPred.Cells.Item(8, 4).Value = 18
$wsPred.Cells.Item(8, 5).Value = -4
$wsPred.Cells.Item(8, 6).Value = -5
$wsPred.Cells.Item(8, 7).Value = "New York Giants +5"
$wsPred.Cells.Item(8, 8).Value = 32
$wsPred.Cells.Item(8, 9).Value = 41.5
$wsPred.Cells.Item(8, 10).Value = "UNDER"

$wsPred.Cells.Item(9, 1).Value = "Denver Broncos"
$wsPred.Cells.Item(9, 2).Value = "Kansas City Chiefs"
$wsPred.Cells.Item(9, 3).Value = 18
$wsPred.Cells.Item(9, 4).Value = 26
$wsPred.Cells.Item(9, 5).Value = -8
$wsPred.Cells.Item(9, 6).Value = -10
$wsPred.Cells.Item(9, 7).Value = "Denver +8"
$wsPred.Cells.Item(9, 8).Value = 44
$wsPred.Cells.Item(9, 9).Value = 47
$wsPred.Cells.Item(9, 10).Value = "UNDER"

$wsPred.Cells.Item(10, 1).Value = "Philadelphia Eagles"
$wsPred.Cells.Item(10, 2).Value = "New York Jets"
$wsPred.Cells.Item(10, 3).Value = 33
$wsPred.Cells.Item(10, 4).Value = 17
$wsPred.Cells.Item(10, 5).Value = 16
$wsPred.Cells.Item(10, 6).Value = 6.5
$wsPred.Cells.Item(10, 7).Value = "Philadelphia -6.5"
$wsPred.Cells.Item(10, 8).Value = 50
$wsPred.Cells.Item(10, 9).Value = 45
$wsPred.Cells.Item(10, 10).Value = "OVER"

$wsPred.Cells.Item(11, 1).Value = "Washington Football Team"
$wsPred.Cells.Item(11, 2).Value = "Las Vegas Raiders"
$wsPred.Cells.Item(11, 3).Value = 21
$wsPred.Cells.Item(11, 4).Value = 30
$wsPred.Cells.Item(11, 5).Value = -9
$wsPred.Cells.Item(11, 6).Value = -2.5
$wsPred.Cells.Item(11, 7).Value = "Las Vegas -2.5"
$wsPred.Cells.Item(11, 8).Value = 51
$wsPred.Cells.Item(11, 9).Value = 49
$wsPred.Cells.Item(11, 10).Value = "OVER"

$wsPred.Cells.Item(12, 1).Value = "Jacksonville Jaguars"
$wsPred.Cells.Item(12, 2).Value = "Los Angeles Rams"
$wsPred.Cells.Item(12, 3).Value = 9
$wsPred.Cells.Item(12, 4).Value = 37
$wsPred.Cells.Item(12, 5).Value = -28
$wsPred.Cells.Item(12, 6).Value = -13
$wsPred.Cells.Item(12, 7).Value = "Los Angeles Rams -13"
$wsPred.Cells.Item(12, 8).Value = 46
$wsPred.Cells.Item(12, 9).Value = 48
$wsPred.Cells.Item(12, 10).Value = "UNDER"

$wsPred.Cells.Item(13, 1).Value = "Baltimore Ravens"
$wsPred.Cells.Item(13, 2).Value = "Pittsburgh Steelers"
$wsPred.Cells.Item(13, 3).Value = 27
$wsPred.Cells.Item(13, 4).Value = 20
$wsPred.Cells.Item(13, 5).Value = 7
$wsPred.Cells.Item(13, 6).Value = 4.5
$wsPred.Cells.Item(13, 7).Value = "Baltimore -4.5"
$wsPred.Cells.Item(13, 8).Value = 47
$wsPred.Cells.Item(13, 9).Value = 44
$wsPred.Cells.Item(13, 10).Value = "OVER"

$wsPred.Cells.Item(14, 1).Value = "San Francisco 49ers"
$wsPred.Cells.Item(14, 2).Value = "Seattle Seahawks"
$wsPred.Cells.Item(14, 3).Value = 25
$wsPred.Cells.Item(14, 4).Value = 16
$wsPred.Cells.Item(14, 5).Value = 9
$wsPred.Cells.Item(14, 6).Value = 3.5
$wsPred.Cells.Item(14, 7).Value = "San Francisco -3.5"
$wsPred.Cells.Item(14, 8).Value = 41
$wsPred.Cells.Item(14, 9).Value = 45.5
$wsPred.Cells.Item(14, 10).Value = "UNDER"

$wsPred.Cells.Item(15, 1).Value = "New England Patriots"
$wsPred.Cells.Item(15, 2).Value = "Buffalo Bills"
$wsPred.Cells.Item(15, 3).Value = 25
$wsPred.Cells.Item(15, 4).Value = 32
$wsPred.Cells.Item(15, 5).Value = -7
$wsPred.Cells.Item(15, 6).Value = -3
$wsPred.Cells.Item(15, 7).Value = "Buffalo -2.5"
$wsPred.Cells.Item(15, 8).Value = 57
$wsPred.Cells.Item(15, 9).Value = 44
$wsPred.Cells.Item(15, 10).Value = "OVER"

# Week 13 only has 14 games (row 16 - the old week-12 15th game - is removed)
$wsPred.Rows.Item(16).Delete()
$wsPred.Range("G15").Select()

# =========================================================
# 2. Results sheet: label existing week-11 rows, append week-12
# =========================================================
$wsRes = $wb.Worksheets.Item("Results")

# Week-11 block (rows 29-43) previously had no value in column A; label it "11"
for ($r = 29; $r -le 43; $r++) {
    $wsRes.Cells.Item($r, 1).NumberFormat = "@"
    $wsRes.Cells.Item($r, 1).Value = "11"
    $wsRes.Cells.Item($r, 1).Style = "Normal"
}

# Week-12 results (rows 44-58), appended below the week-11 block
$wsRes.Cells.Item(44, 1).NumberFormat = "@"
$wsRes.Cells.Item(44, 1).Value = "12"
$wsRes.Cells.Item(44, 1).Style = "Normal"
$wsRes.Cells.Item(44, 2).Value = "Chicago Bears"
$wsRes.Cells.Item(44, 3).Value = "Detroit Lions"
$wsRes.Cells.Item(44, 4).Value = 16
$wsRes.Cells.Item(44, 5).Value = 14
$wsRes.Cells.Item(44, 6).Value = 13
$wsRes.Cells.Item(44, 7).Value = 13
$wsRes.Cells.Item(44, 8).Value = "WIN"
$wsRes.Cells.Item(44, 9).Value = 0
$wsRes.Cells.Item(44, 10).Value = 3
$wsRes.Cells.Item(44, 11).Value = 26
$wsRes.Cells.Item(44, 12).Value = 41.5
$wsRes.Cells.Item(44, 13).Value = 30
$wsRes.Cells.Item(44, 14).Value = "Detroit +3"
$wsRes.Cells.Item(44, 15).Value = "UNDER"
$wsRes.Cells.Item(44, 16).Value = "WIN"

$wsRes.Cells.Item(45, 1).NumberFormat = "@"
$wsRes.Cells.Item(45, 1).Value = "12"
$wsRes.Cells.Item(45, 1).Style = "Normal"
$wsRes.Cells.Item(45, 2).Value = "Las Vegas Raiders"
$wsRes.Cells.Item(45, 3).Value = "Dallas Cowboys"
$wsRes.Cells.Item(45, 4).Value = 36
$wsRes.Cells.Item(45, 5).Value = 33
$wsRes.Cells.Item(45, 6).Value = 20
$wsRes.Cells.Item(45, 7).Value = 40
$wsRes.Cells.Item(45, 8).Value = "LOSS"
$wsRes.Cells.Item(45, 9).Value = -20
$wsRes.Cells.Item(45, 10).Value = -7.5
$wsRes.Cells.Item(45, 11).Value = 60
$wsRes.Cells.Item(45, 12).Value = 51.5
$wsRes.Cells.Item(45, 13).Value = 69
$wsRes.Cells.Item(45, 14).Value = "Dallas -7.5"
$wsRes.Cells.Item(45, 15).Value = "OVER"
$wsRes.Cells.Item(45, 16).Value = "WIN"

$wsRes.Cells.Item(46, 1).NumberFormat = "@"
$wsRes.Cells.Item(46, 1).Value = "12"
$wsRes.Cells.Item(46, 1).Style = "Normal"
$wsRes.Cells.Item(46, 2).Value = "Buffalo Bills"
$wsRes.Cells.Item(46, 3).Value = "New Orleans Saints"
$wsRes.Cells.Item(46, 4).Value = 31
$wsRes.Cells.Item(46, 5).Value = 6
$wsRes.Cells.Item(46, 6).Value = 34
$wsRes.Cells.Item(46, 7).Value = 25
$wsRes.Cells.Item(46, 8).Value = "WIN"
$wsRes.Cells.Item(46, 9).Value = 9
$wsRes.Cells.Item(46, 10).Value = 6
$wsRes.Cells.Item(46, 11).Value = 59
$wsRes.Cells.Item(46, 12).Value = 45.5
$wsRes.Cells.Item(46, 13).Value = 37
$wsRes.Cells.Item(46, 14).Value = "Buffalo -6"
$wsRes.Cells.Item(46, 15).Value = "OVER"
$wsRes.Cells.Item(46, 16).Value = "LOSS"

$wsRes.Cells.Item(47, 1).NumberFormat = "@"
$wsRes.Cells.Item(47, 1).Value = "12"
$wsRes.Cells.Item(47, 1).Style = "Normal"
$wsRes.Cells.Item(47, 2).Value = "Pittsburgh Steelers"
$wsRes.Cells.Item(47, 3).Value = "Cincinnati Bengals"
$wsRes.Cells.Item(47, 4).Value = 10
$wsRes.Cells.Item(47, 5).Value = 41
$wsRes.Cells.Item(47, 6).Value = 18
$wsRes.Cells.Item(47, 7).Value = 34
$wsRes.Cells.Item(47, 8).Value = "WIN"
$wsRes.Cells.Item(47, 9).Value = -16
$wsRes.Cells.Item(47, 10).Value = -4.5
$wsRes.Cells.Item(47, 11).Value = 52
$wsRes.Cells.Item(47, 12).Value = 45
$wsRes.Cells.Item(47, 13).Value = 51
$wsRes.Cells.Item(47, 14).Value = "Cincinnati -4.5"
$wsRes.Cells.Item(47, 15).Value = "OVER"
$wsRes.Cells.Item(47, 16).Value = "WIN"

$wsRes.Cells.Item(48, 1).NumberFormat = "@"
$wsRes.Cells.Item(48, 1).Value = "12"
$wsRes.Cells.Item(48, 1).Style = "Normal"
$wsRes.Cells.Item(48, 2).Value = "New York Jets"
$wsRes.Cells.Item(48, 3).Value = "Houston Texans"
$wsRes.Cells.Item(48, 4).Value = 21
$wsRes.Cells.Item(48, 5).Value = 14
$wsRes.Cells.Item(48, 6).Value = 17
$wsRes.Cells.Item(48, 7).Value = 19
$wsRes.Cells.Item(48, 8).Value = "WIN"
$wsRes.Cells.Item(48, 9).Value = -2
$wsRes.Cells.Item(48, 10).Value = -2.5
$wsRes.Cells.Item(48, 11).Value = 36
$wsRes.Cells.Item(48, 12).Value = 44.5
$wsRes.Cells.Item(48, 13).Value = 35
$wsRes.Cells.Item(48, 14).Value = "New York Jets +2.5"
$wsRes.Cells.Item(48, 15).Value = "UNDER"
$wsRes.Cells.Item(48, 16).Value = "WIN"

$wsRes.Cells.Item(49, 1).NumberFormat = "@"
$wsRes.Cells.Item(49, 1).Value = "12"
$wsRes.Cells.Item(49, 1).Style = "Normal"
$wsRes.Cells.Item(49, 2).Value = "Carolina Panthers"
$wsRes.Cells.Item(49, 3).Value = "Miami Dolphins"
$wsRes.Cells.Item(49, 4).Value = 10
$wsRes.Cells.Item(49, 5).Value = 33
$wsRes.Cells.Item(49, 6).Value = 19
$wsRes.Cells.Item(49, 7).Value = 12
$wsRes.Cells.Item(49, 8).Value = "LOSS"
$wsRes.Cells.Item(49, 9).Value = 7
$wsRes.Cells.Item(49, 10).Value = 2
$wsRes.Cells.Item(49, 11).Value = 31
$wsRes.Cells.Item(49, 12).Value = 42
$wsRes.Cells.Item(49, 13).Value = 43
$wsRes.Cells.Item(49, 14).Value = "Carolina -2"
$wsRes.Cells.Item(49, 15).Value = "UNDER"
$wsRes.Cells.Item(49, 16).Value = "LOSS"

$wsRes.Cells.Item(50, 1).NumberFormat = "@"
$wsRes.Cells.Item(50, 1).Value = "12"
$wsRes.Cells.Item(50, 1).Style = "Normal"
$wsRes.Cells.Item(50, 2).Value = "Philadelphia Eagles"
$wsRes.Cells.Item(50, 3).Value = "New York Giants"
$wsRes.Cells.Item(50, 4).Value = 7
$wsRes.Cells.Item(50, 5).Value = 13
$wsRes.Cells.Item(50, 6).Value = 30
$wsRes.Cells.Item(50, 7).Value = 18
$wsRes.Cells.Item(50, 8).Value = "LOSS"
$wsRes.Cells.Item(50, 9).Value = 12
$wsRes.Cells.Item(50, 10).Value = 3.5
$wsRes.Cells.Item(50, 11).Value = 48
$wsRes.Cells.Item(50, 12).Value = 45.5
$wsRes.Cells.Item(50, 13).Value = 20
$wsRes.Cells.Item(50, 14).Value = "Philadelphia -3.5"
$wsRes.Cells.Item(50, 15).Value = "OVER"
$wsRes.Cells.Item(50, 16).Value = "LOSS"

$wsRes.Cells.Item(51, 1).NumberFormat = "@"
$wsRes.Cells.Item(51, 1).Value = "12"
$wsRes.Cells.Item(51, 1).Style = "Normal"
$wsRes.Cells.Item(51, 2).Value = "Atlanta Falcons"
$wsRes.Cells.Item(51, 3).Value = "Jacksonville Jaguars"
$wsRes.Cells.Item(51, 4).Value = 21
$wsRes.Cells.Item(51, 5).Value = 14
$wsRes.Cells.Item(51, 6).Value = 16
$wsRes.Cells.Item(51, 7).Value = 17
$wsRes.Cells.Item(51, 8).Value = "LOSS"
$wsRes.Cells.Item(51, 9).Value = -1
$wsRes.Cells.Item(51, 10).Value = 2.5
$wsRes.Cells.Item(51, 11).Value = 33
$wsRes.Cells.Item(51, 12).Value = 46
$wsRes.Cells.Item(51, 13).Value = 35
$wsRes.Cells.Item(51, 14).Value = "Jacksonville +2.5"
$wsRes.Cells.Item(51, 15).Value = "UNDER"
$wsRes.Cells.Item(51, 16).Value = "WIN"

$wsRes.Cells.Item(52, 1).NumberFormat = "@"
$wsRes.Cells.Item(52, 1).Value = "12"
$wsRes.Cells.Item(52, 1).Style = "Normal"
$wsRes.Cells.Item(52, 2).Value = "Tampa Bay Buccaneers"
$wsRes.Cells.Item(52, 3).Value = "Indianapolis Colts"
$wsRes.Cells.Item(52, 4).Value = 38
$wsRes.Cells.Item(52, 5).Value = 31
$wsRes.Cells.Item(52, 6).Value = 37
$wsRes.Cells.Item(52, 7).Value = 35
$wsRes.Cells.Item(52, 8).Value = "LOSS"
$wsRes.Cells.Item(52, 9).Value = 2
$wsRes.Cells.Item(52, 10).Value = 3
$wsRes.Cells.Item(52, 11).Value = 72
$wsRes.Cells.Item(52, 12).Value = 53
$wsRes.Cells.Item(52, 13).Value = 69
$wsRes.Cells.Item(52, 14).Value = "Indianapolis +3"
$wsRes.Cells.Item(52, 15).Value = "OVER"
$wsRes.Cells.Item(52, 16).Value = "WIN"

$wsRes.Cells.Item(53, 1).NumberFormat = "@"
$wsRes.Cells.Item(53, 1).Value = "12"
$wsRes.Cells.Item(53, 1).Style = "Normal"
$wsRes.Cells.Item(53, 2).Value = "Tennessee Titans"
$wsRes.Cells.Item(53, 3).Value = "New England Patriots"
$wsRes.Cells.Item(53, 4).Value = 13
$wsRes.Cells.Item(53, 5).Value = 36
$wsRes.Cells.Item(53, 6).Value = 21
$wsRes.Cells.Item(53, 7).Value = 33
$wsRes.Cells.Item(53, 8).Value = "WIN"
$wsRes.Cells.Item(53, 9).Value = -12
$wsRes.Cells.Item(53, 10).Value = -7
$wsRes.Cells.Item(53, 11).Value = 54
$wsRes.Cells.Item(53, 12).Value = 44
$wsRes.Cells.Item(53, 13).Value = 49
$wsRes.Cells.Item(53, 14).Value = "New England -7"
$wsRes.Cells.Item(53, 15).Value = "OVER"
$wsRes.Cells.Item(53, 16).Value = "WIN"

$wsRes.Cells.Item(54, 1).NumberFormat = "@"
$wsRes.Cells.Item(54, 1).Value = "12"
$wsRes.Cells.Item(54, 1).Style = "Normal"
$wsRes.Cells.Item(54, 2).Value = "Los Angeles Chargers"
$wsRes.Cells.Item(54, 3).Value = "Denver Broncos"
$wsRes.Cells.Item(54, 4).Value = 13
$wsRes.Cells.Item(54, 5).Value = 28
$wsRes.Cells.Item(54, 6).Value = 25
$wsRes.Cells.Item(54, 7).Value = 24
$wsRes.Cells.Item(54, 8).Value = "WIN"
$wsRes.Cells.Item(54, 9).Value = 1
$wsRes.Cells.Item(54, 10).Value = 2.5
$wsRes.Cells.Item(54, 11).Value = 49
$wsRes.Cells.Item(54, 12).Value = 48
$wsRes.Cells.Item(54, 13).Value = 41
$wsRes.Cells.Item(54, 14).Value = "Denver +2.5"
$wsRes.Cells.Item(54, 15).Value = "OVER"
$wsRes.Cells.Item(54, 16).Value = "LOSS"

$wsRes.Cells.Item(55, 1).NumberFormat = "@"
$wsRes.Cells.Item(55, 1).Value = "12"
$wsRes.Cells.Item(55, 1).Style = "Normal"
$wsRes.Cells.Item(55, 2).Value = "Los Angeles Rams"
$wsRes.Cells.Item(55, 3).Value = "Green Bay Packers"
$wsRes.Cells.Item(55, 4).Value = 28
$wsRes.Cells.Item(55, 5).Value = 36
$wsRes.Cells.Item(55, 6).Value = 28
$wsRes.Cells.Item(55, 7).Value = 25
$wsRes.Cells.Item(55, 8).Value = "LOSS"
$wsRes.Cells.Item(55, 9).Value = 3
$wsRes.Cells.Item(55, 10).Value = 1
$wsRes.Cells.Item(55, 11).Value = 53
$wsRes.Cells.Item(55, 12).Value = 47
$wsRes.Cells.Item(55, 13).Value = 64
$wsRes.Cells.Item(55, 14).Value = "Los Angeles Rams -1"
$wsRes.Cells.Item(55, 15).Value = "OVER"
$wsRes.Cells.Item(55, 16).Value = "WIN"

$wsRes.Cells.Item(56, 1).NumberFormat = "@"
$wsRes.Cells.Item(56, 1).Value = "12"
$wsRes.Cells.Item(56, 1).Style = "Normal"
$wsRes.Cells.Item(56, 2).Value = "Minnesota Vikings"
$wsRes.Cells.Item(56, 3).Value = "San Francisco 49ers"
$wsRes.Cells.Item(56, 4).Value = 26
$wsRes.Cells.Item(56, 5).Value = 34
$wsRes.Cells.Item(56, 6).Value = 26
$wsRes.Cells.Item(56, 7).Value = 29
$wsRes.Cells.Item(56, 8).Value = "WIN"
$wsRes.Cells.Item(56, 9).Value = -3
$wsRes.Cells.Item(56, 10).Value = -3
$wsRes.Cells.Item(56, 11).Value = 55
$wsRes.Cells.Item(56, 12).Value = 49
$wsRes.Cells.Item(56, 13).Value = 60
$wsRes.Cells.Item(56, 14).Value = "San Francisco -3"
$wsRes.Cells.Item(56, 15).Value = "OVER"
$wsRes.Cells.Item(56, 16).Value = "WIN"

$wsRes.Cells.Item(57, 1).NumberFormat = "@"
$wsRes.Cells.Item(57, 1).Value = "12"
$wsRes.Cells.Item(57, 1).Style = "Normal"
$wsRes.Cells.Item(57, 2).Value = "Cleveland Browns"
$wsRes.Cells.Item(57, 3).Value = "Baltimore Ravens"
$wsRes.Cells.Item(57, 4).Value = 10
$wsRes.Cells.Item(57, 5).Value = 16
$wsRes.Cells.Item(57, 6).Value = 21
$wsRes.Cells.Item(57, 7).Value = 30
$wsRes.Cells.Item(57, 8).Value = "WIN"
$wsRes.Cells.Item(57, 9).Value = -9
$wsRes.Cells.Item(57, 10).Value = -3.5
$wsRes.Cells.Item(57, 11).Value = 51
$wsRes.Cells.Item(57, 12).Value = 46
$wsRes.Cells.Item(57, 13).Value = 26
$wsRes.Cells.Item(57, 14).Value = "Baltimore -3.5"
$wsRes.Cells.Item(57, 15).Value = "OVER"
$wsRes.Cells.Item(57, 16).Value = "LOSS"

$wsRes.Cells.Item(58, 1).NumberFormat = "@"
$wsRes.Cells.Item(58, 1).Value = "12"
$wsRes.Cells.Item(58, 1).Style = "Normal"
$wsRes.Cells.Item(58, 2).Value = "Seattle Seahawks"
$wsRes.Cells.Item(58, 3).Value = "Washington Football Team"
$wsRes.Cells.Item(58, 4).Value = 15
$wsRes.Cells.Item(58, 5).Value = 17
$wsRes.Cells.Item(58, 6).Value = 19
$wsRes.Cells.Item(58, 7).Value = 20
$wsRes.Cells.Item(58, 8).Value = "WIN"
$wsRes.Cells.Item(58, 9).Value = -1
$wsRes.Cells.Item(58, 10).Value = -1
$wsRes.Cells.Item(58, 11).Value = 39
$wsRes.Cells.Item(58, 12).Value = 46.5
$wsRes.Cells.Item(58, 13).Value = 32
$wsRes.Cells.Item(58, 14).Value = "Washington -1"
$wsRes.Cells.Item(58, 15).Value = "UNDER"
$wsRes.Cells.Item(58, 16).Value = "WIN"

$wsRes.Range("P44:P58").Select()

# =========================================================
# 3. WL Record sheet: formulas auto-recalculate; just refresh selection
# =========================================================
$wsWL = $wb.Worksheets.Item("WL Record")
$wsWL.Range("B7").Select()

$wsPred.Select()